# Apply view-count (column F) updates to the "展览" and "全部类型" sheets.
# These mirror the same events (matched by the id= in column H), so both
# sheets receive the same incremented values.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1): row -> new F value
$exhibitionUpdates = @{
    5  = 1416
    9  = 441
    10 = 6623
    14 = 4884
    17 = 6044
    18 = 7789
    20 = 1101
    22 = 4121
    24 = 61
    28 = 1092
    31 = 741
    34 = 2006
    36 = 1295
    40 = 2729
    43 = 5
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (sheet4): row -> new F value
$allTypesUpdates = @{
    10 = 1416
    14 = 441
    15 = 6623
    19 = 4884
    20 = 6044
    21 = 6044
    24 = 4121
    30 = 741
    33 = 2006
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
